$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adiciona a nova coluna M: "Bandeira do Brasil"
$ws.Range("M1").Value = "Bandeira do Brasil"
$ws.Range("M2").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("M4").Value = 0
